$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 10 (shifts existing rows 10-14 down to 11-15)
$ws.Rows.Item(10).Insert()

# Fill the new row 10 with the latest weekly record (copy of former row 10's
# static fields, with the date/price columns updated)
$ws.Cells.Item(10, 1).Value = 9
$ws.Cells.Item(10, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(10, 3).Value = "Metropolitana"
$ws.Cells.Item(10, 4).Value = 44438
$ws.Cells.Item(10, 5).Value = 13
$ws.Cells.Item(10, 6).Value = 100112010
$ws.Cells.Item(10, 7).Value = "Achicoria"
$ws.Cells.Item(10, 8).Value = "Sin especificar"
$ws.Cells.Item(10, 9).Value = "Primera"
$ws.Cells.Item(10, 10).Value = 34
$ws.Cells.Item(10, 11).Value = 5000
$ws.Cells.Item(10, 12).Value = 6000
$ws.Cells.Item(10, 13).Value = 5500
$ws.Cells.Item(10, 14).Value = "`$/caja 16 unidades"
$ws.Cells.Item(10, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(10, 16).Value = 344
$ws.Cells.Item(10, 17).Value = 16
$ws.Cells.Item(10, 18).Value = "Hortaliza"
